$d = $word.ActiveDocument

# Locate the paragraph that contains the "Part 2 of this assignment..." text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Part 2 of this assignment involves*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Part 2 of this assignment...' paragraph"
}

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$rsquo = [char]0x2019

$newXml = "<w:p $wNs>" +
    "<w:r><w:t>Part 2 of this assignment involves implementing the basic functionality of assignment 1 and integrating it into a self-contained GUI. As well as allowing the user to create and modify the layout of roads and traffic lights and control the spawn rate of cars.</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:r><w:t xml:space='preserve'>I have had significant trouble regarding the functionality required for assignment 2, which still have </w:t></w:r>" +
    "<w:r><w:t>not</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> and probabl</w:t></w:r>" +
    "<w:r><w:t>y</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> will not be satisfied. </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>Functionality has been implemented to read and store road and traffic light information from/ to a csv file, but this could still be improved to encode the position on screen that these roads would be drawn. As the program is now, it iterates over each line in the csv, drawing roads consecutively instead of associating positional data of the roads in the csv. </w:t></w:r>" +
    "<w:r><w:t>There is no functionality to allow a user to create or edit a map and the data is currently only able to be read/ stored in one file. Finally, the program still outputs through the console and is unable to update the GUI element, there is significant bottom layer code for this functionality (the repaint() function is quite extensive but is only able to be called during its initialiser and doesn${rsquo}t update with the simulation).</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
    "</w:p>"

$target.Range.InsertXML($newXml)
